$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Target cluster) changes from "Inflammatory-Mac" to "Resolving-Mac" for rows 2-6
$ws.Range("D2:D6").Value = "Resolving-Mac"

# Row 2
$ws.Range("G2").Value = 0.1433473333333333
$ws.Range("H2").Value = 0.430042
$ws.Range("I2").Value = 0.03250327325123204
$ws.Range("J2").Value = 0.03250327325123203
$ws.Range("M2").Value = 0.07218766666666666
$ws.Range("N2").Value = 0.216563
$ws.Range("Q2").Value = 0.01034790951622222
$ws.Range("R2").Value = 0.09313118564599999
$ws.Range("S2").Value = 0.03250327325123204
$ws.Range("T2").Value = 0.03250327325123203

# Row 3
$ws.Range("I3").Value = 0.5075365083813599
$ws.Range("J3").Value = 0.5075365083813598
$ws.Range("M3").Value = 0.07218766666666666
$ws.Range("N3").Value = 0.216563
$ws.Range("Q3").Value = 0.1615819374348889
$ws.Range("R3").Value = 1.454237436914
$ws.Range("S3").Value = 0.5075365083813599
$ws.Range("T3").Value = 0.5075365083813598

# Row 4
$ws.Range("G4").Value = 1.598977333333333
$ws.Range("H4").Value = 4.796932
$ws.Range("I4").Value = 0.362559916388583
$ws.Range("J4").Value = 0.362559916388583
$ws.Range("M4").Value = 0.07218766666666666
$ws.Range("N4").Value = 0.216563
$ws.Range("Q4").Value = 0.1154264427462222
$ws.Range("R4").Value = 1.038837984716
$ws.Range("S4").Value = 0.362559916388583
$ws.Range("T4").Value = 0.362559916388583

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.06293633333333333
$ws.Range("H5").Value = 0.188809
$ws.Range("I5").Value = 0.01427049106666761
$ws.Range("J5").Value = 0.0142704910666676
$ws.Range("M5").Value = 0.07218766666666666
$ws.Range("N5").Value = 0.216563
$ws.Range("Q5").Value = 0.004543227051888889
$ws.Range("R5").Value = 0.040889043467
$ws.Range("S5").Value = 0.01427049106666761
$ws.Range("T5").Value = 0.0142704910666676

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3666226666666666
$ws.Range("H6").Value = 1.099868
$ws.Range("I6").Value = 0.0831298109121576
$ws.Range("J6").Value = 0.08312981091215758
$ws.Range("M6").Value = 0.07218766666666666
$ws.Range("N6").Value = 0.216563
$ws.Range("Q6").Value = 0.02646563485377777
$ws.Range("R6").Value = 0.238190713684
$ws.Range("S6").Value = 0.0831298109121576
$ws.Range("T6").Value = 0.08312981091215758
